# Apply the workbook edits described by the commit:
# "updated docs and default task parameters"
#
# - updates default task-parameter values in row 2 (J2, L2, M2, N2, P2)
# - moves the sheet's visible window / selection (topLeftCell, activeCell)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated default task parameters (row 2) ---
$ws.Range("J2").Value = 3
$ws.Range("L2").Value = 52
$ws.Range("M2").Value = 17
$ws.Range("N2").Value = 6
$ws.Range("P2").Value = 3

# --- View / selection state ---
# Scroll the visible pane so column I is the left-most visible column
# (mirrors topLeftCell moving from K1 to I1).
$win = $excel.ActiveWindow
$win.ScrollColumn = 9
$win.ScrollRow = 1

# Move the active selection from N2 to P3.
$ws.Range("P3").Select()
